$d = $word.ActiveDocument

# Update the base Jira link used in the task-description template text.
$d.Content.Find.Execute(
    "https://jira.edev.pro/browse/AFO-",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://afo.atlassian.net/browse/AFO-",
    2
)

# Minor column-width rebalancing in the main order table (935->936, 1069->1068 dxa).
# Updating Column.Width (in points; 20 dxa = 1 pt) keeps the tblGrid and every
# spanned tc's tcW in sync automatically.
$t = $d.Tables.Item(1)
$t.Columns.Item(4).Width = 46.8
$t.Columns.Item(8).Width = 53.4
